$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value2 = 2627.7144
$ws.Range("I125").Value2 = 397
$ws.Range("K125").Value2 = 3573
$ws.Range("M125").Value2 = -1113

$ws.Range("H132").Value2 = 6016.794
$ws.Range("I132").Value2 = 6095.645
$ws.Range("J132").Value2 = 5202
$ws.Range("K132").Value2 = 18286.935
$ws.Range("L132").Value2 = 15606
$ws.Range("M132").Value2 = -15756.935
$ws.Range("N132").Value2 = -20666

$ws.Range("H138").Value2 = 440796.72
$ws.Range("I138").Value2 = 883.4091
$ws.Range("J138").Value2 = 566486.25
$ws.Range("K138").Value2 = 2650.2273
$ws.Range("L138").Value2 = 1699458.75
$ws.Range("M138").Value2 = 2489.7727
$ws.Range("N138").Value2 = -1709738.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 760.7222
$ws.Range("I2").Value2 = 675.61536
$ws.Range("J2").Value2 = 982
$ws.Range("K2").Value2 = 675.61536
$ws.Range("L2").Value2 = 982
$ws.Range("M2").Value2 = -562.61536
$ws.Range("N2").Value2 = -1208

$ws.Range("H23").Value2 = 74169.836

$ws.Range("H32").Value2 = 5607.617
$ws.Range("I32").Value2 = 5616.8667
$ws.Range("J32").Value2 = 5399.5
$ws.Range("K32").Value2 = 5616.8667
$ws.Range("L32").Value2 = 5399.5
$ws.Range("M32").Value2 = -5329.8667
$ws.Range("N32").Value2 = -5973.5

$ws.Range("H37").Value2 = 24360
$ws.Range("J37").Value2 = 28000
$ws.Range("L37").Value2 = 28000
$ws.Range("N37").Value2 = -28546

$ws.Range("H44").Value2 = 23297.6
$ws.Range("I44").Value2 = 13544
$ws.Range("K44").Value2 = 13544
$ws.Range("M44").Value2 = -13056

$ws.Range("H55").Value2 = 37966.332
$ws.Range("J55").Value2 = 37966.332
$ws.Range("L55").Value2 = 37966.332
$ws.Range("N55").Value2 = -38596.332

$ws.Range("H74").Value2 = 1973.5555
$ws.Range("I74").Value2 = 1443.6666
$ws.Range("J74").Value2 = 3033.3333
$ws.Range("K74").Value2 = 1443.6666
$ws.Range("L74").Value2 = 3033.3333
$ws.Range("M74").Value2 = -569.6666
$ws.Range("N74").Value2 = -4781.3333

$ws.Range("H77").Value2 = 1973.5555
$ws.Range("I77").Value2 = 1443.6666
$ws.Range("J77").Value2 = 3033.3333
$ws.Range("K77").Value2 = 7218.333000000001
$ws.Range("L77").Value2 = 15166.6665
$ws.Range("M77").Value2 = -2850.333000000001
$ws.Range("N77").Value2 = -23902.6665

$ws.Range("H116").Value2 = 760.7222
$ws.Range("I116").Value2 = 675.61536
$ws.Range("J116").Value2 = 982
$ws.Range("K116").Value2 = 675.61536
$ws.Range("L116").Value2 = 982
$ws.Range("M116").Value2 = 1618.38464
$ws.Range("N116").Value2 = -5570

$ws.Range("H122").Value2 = 1367.5264
$ws.Range("I122").Value2 = 1187.2307
$ws.Range("J122").Value2 = 1758.1666
$ws.Range("K122").Value2 = 3561.6921
$ws.Range("L122").Value2 = 5274.4998
$ws.Range("M122").Value2 = -1111.6921
$ws.Range("N122").Value2 = -10174.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 760.7222
$ws.Range("I3").Value2 = 675.61536
$ws.Range("J3").Value2 = 982
$ws.Range("K3").Value2 = 675.61536
$ws.Range("L3").Value2 = 982
$ws.Range("M3").Value2 = -561.61536
$ws.Range("N3").Value2 = -1210

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1720.3103
$ws.Range("I31").Value2 = 2257
$ws.Range("J31").Value2 = 1437.8422
$ws.Range("K31").Value2 = 2257
$ws.Range("L31").Value2 = 1437.8422
$ws.Range("M31").Value2 = -1962
$ws.Range("N31").Value2 = -2027.8422

$ws.Range("H34").Value2 = 1720.3103
$ws.Range("I34").Value2 = 2257
$ws.Range("J34").Value2 = 1437.8422
$ws.Range("K34").Value2 = 2257
$ws.Range("L34").Value2 = 1437.8422
$ws.Range("M34").Value2 = -2055
$ws.Range("N34").Value2 = -1841.8422

$ws.Range("H58").Value2 = 1451
$ws.Range("I58").Value2 = 1138.6471
$ws.Range("K58").Value2 = 1138.6471
$ws.Range("M58").Value2 = -935.6470999999999

$ws.Range("H136").Value2 = 1451
$ws.Range("I136").Value2 = 1138.6471
$ws.Range("K136").Value2 = 3415.9413
$ws.Range("M136").Value2 = -865.9412999999995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value2 = 1130.4706
$ws.Range("J122").Value2 = 1148
$ws.Range("L122").Value2 = 10332
$ws.Range("N122").Value2 = -15232

$ws.Range("H139").Value2 = 1516.1351
$ws.Range("I139").Value2 = 1405.1305
$ws.Range("J139").Value2 = 1698.5
$ws.Range("K139").Value2 = 4215.3915
$ws.Range("L139").Value2 = 5095.5
$ws.Range("M139").Value2 = 924.6085000000003
$ws.Range("N139").Value2 = -15375.5

$ws.Range("H140").Value2 = 26300.395
$ws.Range("I140").Value2 = 79074.30499999999
$ws.Range("K140").Value2 = 237222.915
$ws.Range("M140").Value2 = -232042.915

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 5192.857
$ws.Range("I80").Value2 = 4311.1113
$ws.Range("J80").Value2 = 6780
$ws.Range("K80").Value2 = 4311.1113
$ws.Range("L80").Value2 = 6780
$ws.Range("M80").Value2 = -3313.1113
$ws.Range("N80").Value2 = -8776

$ws.Range("H83").Value2 = 5192.857
$ws.Range("I83").Value2 = 4311.1113
$ws.Range("J83").Value2 = 6780
$ws.Range("K83").Value2 = 21555.5565
$ws.Range("L83").Value2 = 33900
$ws.Range("M83").Value2 = -16563.5565
$ws.Range("N83").Value2 = -43884

$ws.Range("H122").Value2 = 1687.2
$ws.Range("I122").Value2 = 1664.08
$ws.Range("J122").Value2 = 1745
$ws.Range("K122").Value2 = 4992.24
$ws.Range("L122").Value2 = 5235
$ws.Range("M122").Value2 = -2542.24
$ws.Range("N122").Value2 = -10135

$ws.Range("H126").Value2 = 1923.4
$ws.Range("I126").Value2 = 1693
$ws.Range("K126").Value2 = 5079
$ws.Range("M126").Value2 = -2609

$ws.Range("H132").Value2 = 3743.1765
$ws.Range("I132").Value2 = 3942.2
$ws.Range("J132").Value2 = 3458.8572
$ws.Range("K132").Value2 = 11826.6
$ws.Range("L132").Value2 = 10376.5716
$ws.Range("M132").Value2 = -9296.599999999999
$ws.Range("N132").Value2 = -15436.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 1570
$ws.Range("I7").Value2 = 1442.5
$ws.Range("J7").Value2 = 1952.5
$ws.Range("K7").Value2 = 1442.5
$ws.Range("L7").Value2 = 1952.5
$ws.Range("M7").Value2 = -1330.5
$ws.Range("N7").Value2 = -2176.5

$ws.Range("H40").Value2 = 2737.923
$ws.Range("I40").Value2 = 1817.3182
$ws.Range("K40").Value2 = 1817.3182
$ws.Range("M40").Value2 = -1681.3182

$ws.Range("H68").Value2 = 2033
$ws.Range("J68").Value2 = 2000
$ws.Range("L68").Value2 = 2000
$ws.Range("N68").Value2 = -3498

$ws.Range("H71").Value2 = 2033
$ws.Range("J71").Value2 = 2000
$ws.Range("L71").Value2 = 10000
$ws.Range("N71").Value2 = -17488

$ws.Range("H126").Value2 = 1570
$ws.Range("I126").Value2 = 1442.5
$ws.Range("J126").Value2 = 1952.5
$ws.Range("K126").Value2 = 4327.5
$ws.Range("L126").Value2 = 5857.5
$ws.Range("M126").Value2 = -1857.5
$ws.Range("N126").Value2 = -10797.5

$ws.Range("H132").Value2 = 2563.2856
$ws.Range("I132").Value2 = 2210.2354
$ws.Range("K132").Value2 = 6630.706200000001
$ws.Range("M132").Value2 = -4100.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value2 = 17900
$ws.Range("J64").Value2 = 17900
$ws.Range("L64").Value2 = 17900
$ws.Range("N64").Value2 = -18396

$ws.Range("H67").Value2 = 17900
$ws.Range("J67").Value2 = 17900
$ws.Range("L67").Value2 = 17900
$ws.Range("N67").Value2 = -19616

$ws.Range("H81").Value2 = 243
$ws.Range("I81").Value2 = 243
$ws.Range("K81").Value2 = 486
$ws.Range("M81").Value2 = 575

$ws.Range("H84").Value2 = 243
$ws.Range("I84").Value2 = 243
$ws.Range("K84").Value2 = 2430
$ws.Range("M84").Value2 = 2874
